# Insert a new weekly price record as row 270 in the "Cilantro" sheet.
# This pushes the former rows 270-347 down to 271-348 (dimension becomes
# A1:R348) and populates the new row 270 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 270, shifting existing data down.
$ws.Rows.Item(270).Insert()

# Populate the new row 270 with the new record's data.
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44841
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = 100112040
$ws.Range("G270").Value = "Cilantro"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 240
$ws.Range("K270").Value = 10000
$ws.Range("L270").Value = 10500
$ws.Range("M270").Value = 10250
$ws.Range("N270").Value = "$/caja 36 atados"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 285
$ws.Range("Q270").Value = 36
$ws.Range("R270").Value = "Hortaliza"
